$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "35÷5=7, 0"
$t.Cell(1, 2).Range.Text = "82÷9=9, 1"
$t.Cell(1, 3).Range.Text = "23÷5=4, 3"
$t.Cell(1, 4).Range.Text = "42÷5=8, 2"
$t.Cell(1, 5).Range.Text = "88÷7=12, 4"

$t.Cell(5, 1).Range.Text = "85÷2=42, 1"
$t.Cell(5, 2).Range.Text = "59÷4=14, 3"
$t.Cell(5, 3).Range.Text = "93÷9=10, 3"
$t.Cell(5, 4).Range.Text = "96÷5=19, 1"
$t.Cell(5, 5).Range.Text = "75÷7=10, 5"

$t.Cell(9, 1).Range.Text = "59÷4=14, 3"
$t.Cell(9, 2).Range.Text = "75÷9=8, 3"
$t.Cell(9, 3).Range.Text = "48÷8=6, 0"
$t.Cell(9, 4).Range.Text = "87÷8=10, 7"
$t.Cell(9, 5).Range.Text = "57÷4=14, 1"

$t.Cell(13, 1).Range.Text = "56÷9=6, 2"
$t.Cell(13, 2).Range.Text = "11÷7=1, 4"
$t.Cell(13, 3).Range.Text = "15÷6=2, 3"
$t.Cell(13, 4).Range.Text = "40÷5=8, 0"
$t.Cell(13, 5).Range.Text = "67÷5=13, 2"

$t.Cell(17, 1).Range.Text = "22÷7=3, 1"
$t.Cell(17, 2).Range.Text = "79÷6=13, 1"
$t.Cell(17, 3).Range.Text = "63÷9=7, 0"
$t.Cell(17, 4).Range.Text = "76÷7=10, 6"
$t.Cell(17, 5).Range.Text = "57÷5=11, 2"
